# Apply "Uren bijwerking referentiewaarde toegevoegt" edit to uren DAS.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Week block rows 18-24 ---
# Move/rewrite the note text that used to sit on R20 up onto R19,
# updating its wording at the same time (keeps the shared-string slot).
$ws.Range("R20").Value2 = "Github opzetten, bespreken met Ernst en Jo over data opschoning, verder R programeren"
$ws.Range("R20").Cut($ws.Range("R19"))

# "Notes" header label for this block (style copied from the first week's R3)
$ws.Range("R3").Copy()
$ws.Range("R18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R18").Value2 = "Notes"

# Hours added to the existing week table
$ws.Range("J20").Value2 = 6
$ws.Range("O20").Value2 = 4
$ws.Range("G24").Value2 = 4
$ws.Range("K24").Value2 = 2
$ws.Range("M24").Value2 = 2

# --- Week block rows 26-32 ---
$ws.Range("R3").Copy()
$ws.Range("R26").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R26").Value2 = "Notes"

$ws.Range("R6").Copy()
$ws.Range("R27").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R27").Value2 = "R Markdown programmeren, verslag verwerken"

$ws.Range("C28").Value2 = 5
$ws.Range("H28").Value2 = 2
$ws.Range("J28").Value2 = 2
$ws.Range("K28").Value2 = 5

$ws.Range("D31").Value2 = 2

$ws.Range("C32").Value2 = 2
$ws.Range("F32").Value2 = 2
$ws.Range("G32").Value2 = 4
$ws.Range("H32").Value2 = 1
$ws.Range("I32").Value2 = 2
$ws.Range("L32").Value2 = 4

# Update selection to match the saved file's cursor position
$ws.Range("K13").Select()
